$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Professionalism")

# Remove the old "Work Process" rows (7-15) entirely, including their
# row-level formatting (e.g. the ht=30 on the old row 9), so we can
# rebuild the block cleanly at the new row numbers.
$ws.Range("A7:A15").EntireRow.Delete()

# Row 7: section header + the "End fo Day" note in column E (unchanged text)
$ws.Range("A7").Value = "Work Process Enumeration"
$ws.Range("A7").Font.Name = "Arial"
$ws.Range("A7").Font.Size = 12
$ws.Range("A7").HorizontalAlignment = -4131
$ws.Range("A7").WrapText = $true
$ws.Range("E7").Value = "End fo Day"

# Rows 8-10: the old combined "Review, Release, Version Control Check
# Lists" item is split into three distinct checklist items
$ws.Range("A8").Value = "Release Check List"
$ws.Range("A8").Font.Name = "Arial"
$ws.Range("A8").Font.Size = 12
$ws.Range("A8").HorizontalAlignment = -4131
$ws.Range("A8").WrapText = $true

$ws.Range("A9").Value = "Version Control Check List"
$ws.Range("A9").Font.Name = "Arial"
$ws.Range("A9").Font.Size = 12
$ws.Range("A9").HorizontalAlignment = -4131
$ws.Range("A9").WrapText = $true

$ws.Range("A10").Value = "Error Mitigation Stratagy"
$ws.Range("A10").Font.Name = "Arial"
$ws.Range("A10").Font.Size = 12
$ws.Range("A10").HorizontalAlignment = -4131
$ws.Range("A10").WrapText = $true

# Row 12: plain (unstyled) entry, same text as before
$ws.Range("A12").Value = "Work Day Enumeration"

# Rows 13-14: "Task List Format" -> "Task List Template", plus a new
# "Responsibility Guidelines" item
$ws.Range("A13").Value = "Task List Template"
$ws.Range("A13").Font.Name = "Arial"
$ws.Range("A13").Font.Size = 12
$ws.Range("A13").HorizontalAlignment = -4131
$ws.Range("A13").WrapText = $true

$ws.Range("A14").Value = "Responsibility Guidelines"
$ws.Range("A14").Font.Name = "Arial"
$ws.Range("A14").Font.Size = 12
$ws.Range("A14").HorizontalAlignment = -4131
$ws.Range("A14").WrapText = $true

# Rows 16-17: new checklist / guideline items
$ws.Range("A16").Value = "Review Check List"
$ws.Range("A16").Font.Name = "Arial"
$ws.Range("A16").Font.Size = 12
$ws.Range("A16").HorizontalAlignment = -4131
$ws.Range("A16").WrapText = $true

$ws.Range("A17").Value = "Conflict Guidelines"
$ws.Range("A17").Font.Name = "Arial"
$ws.Range("A17").Font.Size = 12
$ws.Range("A17").HorizontalAlignment = -4131
$ws.Range("A17").WrapText = $true

# Row 19: plain (unstyled) entry, same text as before
$ws.Range("A19").Value = "Project Scoping"

# Update the sheet selection to match the new block
$ws.Range("A7:A19").Select()
